# feat: add 2022-Q1 data
#
# Insert a new "2022-Q1" quarterly sheet (with the same layout as the other
# per-quarter fund sheets) right before the "总计" (total) summary sheet, and
# prepend a corresponding new row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# Duplicate the most recent quarterly sheet (brings along sheetPr/pageMargins/
# styles unchanged) and drop it in right before "总计", then rename it.
$templateSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# ---- Header row ----
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- Data rows ----
# Row 2: 002423
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2:B2, D2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "002423"
$newSheet.Range("C2").Value = "华宝兴业标普美国消费(QDII-LOF)美元"
$newSheet.Range("D2").Value = "3.62"
$newSheet.Range("E2").Value = "94.37"
$newSheet.Range("F2").Value = "4.22"
$newSheet.Range("G2").Value = "0.1528"
$newSheet.Range("H2").Value = 4
$newSheet.Range("B2:B2, D2:G2").ClearFormats()

# Row 3: 162415
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3:B3, D3:G3").NumberFormat = "@"
$newSheet.Range("B3").Value = "162415"
$newSheet.Range("C3").Value = "华宝标普美国消费(QDII-LOF)人民币A"
$newSheet.Range("D3").Value = "3.62"
$newSheet.Range("E3").Value = "94.37"
$newSheet.Range("F3").Value = "4.22"
$newSheet.Range("G3").Value = "0.1528"
$newSheet.Range("H3").Value = 4
$newSheet.Range("B3:B3, D3:G3").ClearFormats()

# Row 4: 009975
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4:B4, D4:G4").NumberFormat = "@"
$newSheet.Range("B4").Value = "009975"
$newSheet.Range("C4").Value = "华宝标普美国消费(QDII-LOF)人民币C"
$newSheet.Range("D4").Value = "0.61"
$newSheet.Range("E4").Value = "94.37"
$newSheet.Range("F4").Value = "4.22"
$newSheet.Range("G4").Value = "0.0257"
$newSheet.Range("H4").Value = 4
$newSheet.Range("B4:B4, D4:G4").ClearFormats()

# ---- Update the "总计" (total) summary sheet: insert a new top data row ----
# Re-fetch the sheet by name: inserting the new sheet above shifted sheet
# positions out from under the earlier $totalSheet handle.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows("2:2").Insert()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.33

# Re-sequence the index column (A) for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# Restore the originally-active sheet/selection.
$wb.Worksheets.Item("2020-Q4").Select()
